$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.3491903333333333
$ws.Range("H2").Value = 1.047571
$ws.Range("I2").Value = 0.008130334326258625
$ws.Range("J2").Value = 0.008130334326258625
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 26.91891634927444
$ws.Range("R2").Value = 242.27024714347
$ws.Range("S2").Value = 0.001954383321643403
$ws.Range("T2").Value = 0.001954383321643403

$ws.Range("G3").Value = 0.3491903333333333
$ws.Range("H3").Value = 1.047571
$ws.Range("I3").Value = 0.008130334326258625
$ws.Range("J3").Value = 0.008130334326258625
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 35.47076709643911
$ws.Range("R3").Value = 319.236903867952
$ws.Range("S3").Value = 0.002575269922447926
$ws.Range("T3").Value = 0.002575269922447925

$ws.Range("G4").Value = 0.3491903333333333
$ws.Range("H4").Value = 1.047571
$ws.Range("I4").Value = 0.008130334326258625
$ws.Range("J4").Value = 0.008130334326258625
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 49.59438190956978
$ws.Range("R4").Value = 446.3494371861281
$ws.Range("S4").Value = 0.003600681082167297
$ws.Range("T4").Value = 0.003600681082167297

$ws.Range("I5").Value = 0.801301577139928
$ws.Range("J5").Value = 0.8013015771399279
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 2653.048356929924
$ws.Range("R5").Value = 23877.43521236931
$ws.Range("S5").Value = 0.192618209181256
$ws.Range("T5").Value = 0.192618209181256

$ws.Range("I6").Value = 0.801301577139928
$ws.Range("J6").Value = 0.8013015771399279
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.2538109464642574
$ws.Range("T6").Value = 0.2538109464642574

$ws.Range("I7").Value = 0.801301577139928
$ws.Range("J7").Value = 0.8013015771399279
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.3548724214944146
$ws.Range("T7").Value = 0.3548724214944146

$ws.Range("I8").Value = 0.1905680885338134
$ws.Range("J8").Value = 0.1905680885338134
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 630.95639468537
$ws.Range("R8").Value = 5678.607552168331
$ws.Range("S8").Value = 0.04580907486978303
$ws.Range("T8").Value = 0.04580907486978303

$ws.Range("I9").Value = 0.1905680885338134
$ws.Range("J9").Value = 0.1905680885338134
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("S9").Value = 0.06036212619135434
$ws.Range("T9").Value = 0.06036212619135434

$ws.Range("I10").Value = 0.1905680885338134
$ws.Range("J10").Value = 0.1905680885338134
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("S10").Value = 0.08439688747267604
$ws.Range("T10").Value = 0.08439688747267604

